$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New fitted peak values (Peak_Cent_SO2, Peak_Area_SO2, Peak_Height_SO2)
# for rows 2-10 (data rows 4-12), reflecting updated comments/fits.
$data = @(
    @(2, 1151.160126426616, 55.5762277319274,  26.19545941933148),
    @(3, 1150.607814272824, 75.65668064408545, 54.37831380276303),
    @(4, 1150.65918054404,  645.470942788454,  432.7812345022195),
    @(5, 1150.648215614395, 359.0363137542182, 233.9499953839853),
    @(6, 1150.622331832785, 1121.728413380102, 727.7326777802815),
    @(7, 1150.620833588373, 781.5099153585764, 496.1253298637891),
    @(8, 1150.601837921057, 622.6045691208002, 396.5477527483168),
    @(9, 1150.581569008737, 1114.461012785161, 717.1303478810305),
    @(10,1150.601624386216, 116.5589603381525, 73.53203131267618)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
